$d = $word.ActiveDocument

# --- Locate the "Accomplishment" list item that needs to change -----------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*1. Prepare paper for ASP-DAC*") {
        $target = $p
    }
}

if ($target -eq $null) {
    Write-Host "ERROR: could not find the 'Prepare paper for ASP-DAC' paragraph"
} else {
    # Replace its text in place (keeps the paragraph's own formatting).
    $target.Range.Find.Execute(
        "1. Prepare paper for ASP-DAC (DDL:2019.7.12);", $true, $false, $false,
        $false, $false, $true, 1, $false,
        "1. Get familiar with Xilinx Vivado Framework", 2)
}

# --- Re-acquire the (now retitled) paragraph and append two new items -----
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*1. Get familiar with Xilinx Vivado Framework*") {
        $target = $p
    }
}

# Item 2: "2. Re-implement Xilinx Soc Labs"
$target.Range.InsertParagraphAfter()
$item2 = $target.Next()
$item2.Range.InsertBefore("2. Re-implement Xilinx Soc Labs")

# Item 3: "3. XNORAM patent writing"
$item2.Range.InsertParagraphAfter()
$item3 = $item2.Next()
$item3.Range.InsertBefore("3. XNORAM patent writing")
